# Daily attendance processing - 2026-01-10 17:31:56
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (column G) wherever it currently reads
# "System, dnasr281@gmail.com" so it reads "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G is the 7th column ("Recorded By")
for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
